# Insert a new row before row 177, shifting existing rows 177..289 down to 178..290,
# then populate the newly inserted row 177 with the new data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 177 (this shifts row 177 -> 178, ..., 289 -> 290)
$ws.Rows.Item(177).Insert()

# Populate the new row 177 with the new record's values
$ws.Cells.Item(177, 1).Value = 3
$ws.Cells.Item(177, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(177, 3).Value = "Coquimbo"
$ws.Cells.Item(177, 4).Value = 44596
$ws.Cells.Item(177, 5).Value = 5
$ws.Cells.Item(177, 6).Value = 100114013
$ws.Cells.Item(177, 7).Value = "Zanahoria"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 110
$ws.Cells.Item(177, 11).Value = 7000
$ws.Cells.Item(177, 12).Value = 7000
$ws.Cells.Item(177, 13).Value = 7000
$ws.Cells.Item(177, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(177, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(177, 16).Value = 350
$ws.Cells.Item(177, 17).Value = 20
$ws.Cells.Item(177, 18).Value = "Hortaliza"
